$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "49.781.47"
Set-TextValue "E2" "  -0.55%  "

Set-TextValue "D3" "2.652.22"
Set-TextValue "E3" "  +0.14%  "

Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.03%  "

Set-TextValue "D5" "113.05"
Set-TextValue "E5" "  -0.88%  "

Set-TextValue "D6" "327.45"

Set-TextValue "E7" "  -1.10%  "

Set-TextValue "E8" "  -0.09%  "

Set-TextValue "D9" "0.552"
Set-TextValue "E9" "  -0.94%  "

Set-TextValue "D10" "39.85"
Set-TextValue "E10" "  -2.69%  "

Set-TextValue "D11" "20.05"
Set-TextValue "E11" "  -0.51%  "

Set-TextValue "D12" "0.0818"
Set-TextValue "E12" "  -0.64%  "

Set-TextValue "E13" "  +2.18%  "

Set-TextValue "D14" "7.59"
Set-TextValue "E14" "  +2.68%  "

Set-TextValue "D15" "3.064.91"
Set-TextValue "E15" "  +0.06%  "

Set-TextValue "D16" "2.654.90"
Set-TextValue "E16" "  +0.27%  "

Set-TextValue "D17" "0.860"
Set-TextValue "E17" "  -1.54%  "

Set-TextValue "D18" "49.713.33"
Set-TextValue "E18" "  -0.53%  "

Set-TextValue "D19" "13.37"
Set-TextValue "E19" "  +0.99%  "

Set-TextValue "D20" "6.71"
Set-TextValue "E20" "  -1.19%  "

Set-TextValue "E21" "  -0.13%  "

Set-TextValue "D22" "0.0₃0951"
Set-TextValue "E22" "  -0.63%  "

Set-TextValue "D23" "269.47"
Set-TextValue "E23" "  -2.42%  "

Set-TextValue "D24" "69.23"
Set-TextValue "E24" "  -4.20%  "

Set-TextValue "E25" "  -0.57%  "

Set-TextValue "D26" "26.22"
Set-TextValue "E26" "  -2.35%  "

Set-TextValue "D27" "1.00"
Set-TextValue "E27" "  +0.03%  "

Set-TextValue "D28" "10.23"
Set-TextValue "E28" "  +1.71%  "

Set-TextValue "E29" "  -0.78%  "

Set-TextValue "B30" "InjectiveProtocol"
Set-TextValue "C30" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D30" "34.97"
Set-TextValue "E30" "  -3.99%  "

Set-TextValue "B31" "Kaspa"
Set-TextValue "C31" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D31" "0.138"
Set-TextValue "E31" "  -2.59%  "

Set-TextValue "D32" "49.60"
Set-TextValue "E32" "  -1.10%  "

Set-TextValue "E33" "  +0.26%  "

Set-TextValue "E34" "  +0.72%  "

Set-TextValue "D35" "19.24"
Set-TextValue "E35" "  -1.77%  "

Set-TextValue "E36" "  -0.20%  "

Set-TextValue "E37" "  -1.71%  "

Set-TextValue "E38" "  -1.43%  "

Set-TextValue "D39" "3.14"
Set-TextValue "E39" "  +1.24%  "

Set-TextValue "B40" "EnergySwap"
Set-TextValue "C40" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D40" "23.73"
Set-TextValue "E40" "  +8.02%  "

Set-TextValue "B41" "Monero"
Set-TextValue "C41" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D41" "128.93"
Set-TextValue "E41" "  +4.38%  "

Set-TextValue "D42" "0.0346"
Set-TextValue "E42" "  +9.26%  "

Set-TextValue "D43" "2.28"
Set-TextValue "E43" "  +2.32%  "

Set-TextValue "E44" "  -0.62%  "

Set-TextValue "D45" "3.33"
Set-TextValue "E45" "  -0.34%  "

Set-TextValue "D46" "2.067.31"
Set-TextValue "E46" "  -0.86%  "

Set-TextValue "D47" "2.12"
Set-TextValue "E47" "  +6.42%  "

Set-TextValue "E48" "  -2.06%  "

Set-TextValue "E49" "  -2.32%  "

Set-TextValue "D50" "5.28"
Set-TextValue "E50" "  -2.05%  "

Set-TextValue "D51" "59.11"
Set-TextValue "E51" "  -1.13%  "
